$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.036.64"
$ws.Range("E2").Value = "  +2.46%  "

$ws.Range("D3").Value = "3.040.31"
$ws.Range("E3").Value = "  +1.49%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'594.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").Value = "'154.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.66%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.038.06"
$ws.Range("E8").Value = "  +1.44%  "

$ws.Range("D9").Value = "'0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("D10").Value = "'6.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.68%  "

$ws.Range("E11").Value = "  +3.18%  "

$ws.Range("D12").Value = "'0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.01%  "

$ws.Range("D13").Value = "'0.0000235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.54%  "

$ws.Range("D14").Value = "'35.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.72%  "

$ws.Range("E15").Value = "  +2.01%  "

$ws.Range("D16").Value = "3.543.83"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "63.111.64"
$ws.Range("E17").Value = "  +2.62%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("D19").Value = "3.045.79"
$ws.Range("E19").Value = "  +1.53%  "

$ws.Range("D20").Value = "'454.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").Value = "'14.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.87%  "

$ws.Range("D22").Value = "'0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.81%  "

$ws.Range("D23").Value = "'7.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.44%  "

$ws.Range("D24").Value = "'83.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.91%  "

$ws.Range("D25").Value = "'11.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.31%  "

$ws.Range("D26").Value = "'2.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.19%  "

$ws.Range("D27").Value = "'12.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.32%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").Value = "'7.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.42%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'2.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.36%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").Value = "'27.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("D35").Value = "0.0₃0865"
$ws.Range("E35").Value = "  +4.36%  "

$ws.Range("D36").Value = "'1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.13%  "

$ws.Range("D37").Value = "'5.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.78%  "

$ws.Range("D38").Value = "'3.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.74%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'2.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.19%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.130"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.61%  "

$ws.Range("D41").Value = "'50.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "'9.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "

$ws.Range("E43").Value = "  +11.99%  "

$ws.Range("D44").Value = "'43.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.86%  "

$ws.Range("D45").Value = "'396.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("D46").Value = "'0.0361"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.18%  "

$ws.Range("D47").Value = "2.723.49"
$ws.Range("E47").Value = "  +1.08%  "

$ws.Range("D48").Value = "'132.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.70%  "

$ws.Range("D49").Value = "'2.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.48%  "

$ws.Range("D51").Value = "'24.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.99%  "
